# Daily attendance processing - 2025-10-21 06:28:09
# Normalises the "Recorded By" column (G) so that the actual recorder
# (an email address / named user) is listed before the "System" /
# "system" marker instead of after it. Technically this swaps the first
# two comma-separated entries of every non-empty G-column cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    $first = $parts[0]
    $second = $parts[1]
    $parts[0] = $second
    $parts[1] = $first

    $newText = $parts -join ", "
    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
